$wb = $excel.ActiveWorkbook

# --- "invoice line items" sheet: remove the calculated "line total" column (F) ---
$wsLineItems = $wb.Worksheets.Item("invoice line items")
$wsLineItems.Columns.Item(6).ClearFormats()
$wsLineItems.Columns.Item(6).Delete()
$wsLineItems.Range("F2").Select()

# --- "products" sheet: remove the "Inventory" column (D) ---
$wsProducts = $wb.Worksheets.Item("products")
$wsProducts.Columns.Item(4).ClearFormats()
$wsProducts.Columns.Item(4).Delete()
$wsProducts.Range("D1").Select()

# --- "READ ME" sheet: fix the header/footer font name typo ---
$wsReadMe = $wb.Worksheets.Item("READ ME")
$wsReadMe.PageSetup.CenterHeader = "&""Times New Roman,Regular""&12&A"
$wsReadMe.PageSetup.CenterFooter = "&""Times New Roman,Regular""&12Page &P"
